# Apply the "add monte_carlo and update database" data refresh to the
# Overview sheet of the yearly (rial) income-statement workbook.
#
# The commit replaces a block of placeholder zeros / "-" marks (rows 11-27,
# columns D:H) with the actual reported financial figures, and also turns a
# handful of "-" (no data) text cells into real 0 values now that data is
# known for those periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> value for the five period columns (D..H).
# Row labels (column B), for reference:
#  11 Sales, 12 COGS, 13 Gross profit, 14 G&A expenses,
#  15 Impairment of receivables (exceptional expense),
#  16 Other operating income (expense) net, 17 Operating profit (loss),
#  18 Finance costs, 19 Other non-operating income (expense) net,
#  20 Profit before tax, 21 Tax, 22 Net profit from continuing operations,
#  23 Profit (loss) from discontinued operations, 24 Net profit (loss),
#  25 EPS after tax, 26 Capital, 27 EPS based on latest capital.

$data = @{
    11 = @{ D = 1344450;  E = 2010510;  F = 2535300;  G = 2539337;  H = 5998187 }
    12 = @{ D = -976351;  E = -1514534; F = -1992718; G = -2118070; H = -4526978 }
    13 = @{ D = 368099;   E = 495976;   F = 542582;   G = 421267;   H = 1471209 }
    14 = @{ D = -13335;   E = -48008;   F = -52276;   G = -106372;  H = -178034 }
    15 = @{              E = 0;        F = 0;        G = 0;        H = 0 }
    16 = @{ D = 14294;    E = 6984;     F = 9285;     G = 16764;    H = 14539 }
    17 = @{ D = 369058;   E = 454952;   F = 499591;   G = 331659;   H = 1307714 }
    18 = @{ D = -32505;   E = -25966;   F = -26087;   G = -72263;   H = -142054 }
    19 = @{ D = 32910;    E = 90669;    F = 631791;   G = 1287536;  H = 186999 }
    20 = @{ D = 369463;   E = 519655;   F = 1105295;  G = 1546932;  H = 1352659 }
    21 = @{ D = -43309;   E = -61852;   F = -75069;   G = -62320;   H = 0 }
    22 = @{ D = 326154;   E = 457803;   F = 1030226;  G = 1484612;  H = 1352659 }
    23 = @{ D = 0;        E = 0;        F = 0;        G = 0;        H = 0 }
    24 = @{ D = 326154;   E = 457803;   F = 1030226;  G = 1484612;  H = 1352659 }
    25 = @{ D = 314;      E = 344;      F = 775;      G = 322;      H = 293 }
    26 = @{ D = 1040000;  E = 1330000;  F = 1330000;  G = 4612528;  H = 4612528 }
    27 = @{ D = 53;       E = 75;       F = 169;      G = 243;      H = 222 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
